$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from column P (rows 3-33) into column Q so the new
# "2020" year column visually matches its neighbours, then overwrite
# the values/header with the new data.
$ws.Range("P3:P33").Copy($ws.Range("Q3:Q33"))

$ws.Range("Q3").Value = 2020
$ws.Range("Q4").Value = 53.463696812512026
$ws.Range("Q5").Value = 46.05204738706685
$ws.Range("Q6").Value = 60.76705279190513
$ws.Range("Q7").Value = 46.609654277955656
$ws.Range("Q8").Value = 39.785591828762811
$ws.Range("Q9").Value = 53.693996785869842
$ws.Range("Q10").Value = 49.132459991853935
$ws.Range("Q11").Value = 42.132308166831223
$ws.Range("Q12").Value = 56.225753650646354
$ws.Range("Q13").Value = 28.457427087863305
$ws.Range("Q14").Value = 20.524708126577082
$ws.Range("Q15").Value = 36.325895173845353
$ws.Range("Q16").Value = 37.816151622141014
$ws.Range("Q17").Value = 29.032279844170926
$ws.Range("Q18").Value = 46.928626462141906
$ws.Range("Q19").Value = 51.38232216208695
$ws.Range("Q20").Value = 45.862881450184311
$ws.Range("Q21").Value = 57.0280888993139
$ws.Range("Q22").Value = 44.951834666409091
$ws.Range("Q23").Value = 38.216466887636237
$ws.Range("Q24").Value = 51.83682668469686
$ws.Range("Q25").Value = 82.176148450436926
$ws.Range("Q26").Value = 66.965035434789911
$ws.Range("Q27").Value = 96.931980629894966
$ws.Range("Q28").Value = 56.391242440049062
$ws.Range("Q29").Value = 50.844030930786069
$ws.Range("Q30").Value = 61.300998533028128
$ws.Range("Q31").Value = 54.829571415516767
$ws.Range("Q32").Value = 58.407045187583961
$ws.Range("Q33").Value = 51.452932817170577

# Match the author's final selection/scroll state recorded in the sheet view.
[void]$ws.Range("S34").Select()
